$d = $word.ActiveDocument

# Locate the paragraph that holds the "Обязательства в ПАО «БАНК СГБ»" label
# so the Find/replace below stays scoped to just this paragraph instead of
# touching the whole document story.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Обязательства в ПАО*БАНК СГБ*") {
        $target = $p
        break
    }
}

$r = $target.Range

# Replace the old placeholder run (label + underscores) with the label text
# followed by a trailing space, preserving the original run formatting.
$r.Find.Execute(
    "Обязательства в ПАО «БАНК СГБ»  (рассмотрено, выдано):_________________________________________",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Обязательства в ПАО «БАНК СГБ»  (рассмотрено, выдано): {issue.total_bank_liabilities_vol} руб.",
    2)

# Re-select the merge-field portion we just inserted and underline it,
# which splits it off into its own run with w:u val="single".
$u = $target.Range
$u.Find.Execute(
    "{issue.total_bank_liabilities_vol} руб.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
$u.Font.Underline = 1
